# Auto-generated Excel COM-interop script to apply scheduled-runner value updates
# to the Hades_Profits market-price tracking workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 887.4681
$ws.Range("J17").Value = 887.4681
$ws.Range("L17").Value = 2662.4043
$ws.Range("N17").Value = -2998.4043
$ws.Range("H132").Value = 1169772.8
$ws.Range("I132").Value = 3780.3333
$ws.Range("K132").Value = 11340.9999
$ws.Range("M132").Value = -8810.999899999999
$ws.Range("H138").Value = 1962765.1
$ws.Range("I138").Value = 1240.5577
$ws.Range("J138").Value = 5053652
$ws.Range("K138").Value = 3721.6731
$ws.Range("L138").Value = 15160956
$ws.Range("M138").Value = 1418.3269
$ws.Range("N138").Value = -15171236

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1866.7667
$ws.Range("I2").Value = 1536.3684
$ws.Range("J2").Value = 2437.4546
$ws.Range("K2").Value = 1536.3684
$ws.Range("L2").Value = 2437.4546
$ws.Range("M2").Value = -1423.3684
$ws.Range("N2").Value = -2663.4546
$ws.Range("H32").Value = 1689.25
$ws.Range("I32").Value = 1409.3937
$ws.Range("K32").Value = 1409.3937
$ws.Range("M32").Value = -1122.3937
$ws.Range("H61").Value = 19648536
$ws.Range("I61").Value = 22751156
$ws.Range("K61").Value = 22751156
$ws.Range("M61").Value = -22750944
$ws.Range("H116").Value = 1866.7667
$ws.Range("I116").Value = 1536.3684
$ws.Range("J116").Value = 2437.4546
$ws.Range("K116").Value = 1536.3684
$ws.Range("L116").Value = 2437.4546
$ws.Range("M116").Value = 757.6315999999999
$ws.Range("N116").Value = -7025.4546
$ws.Range("H132").Value = 73353.8
$ws.Range("I132").Value = 50922
$ws.Range("J132").Value = 123202.22
$ws.Range("K132").Value = 152766
$ws.Range("L132").Value = 369606.66
$ws.Range("M132").Value = -150236
$ws.Range("N132").Value = -374666.66
$ws.Range("H136").Value = 19648536
$ws.Range("I136").Value = 22751156
$ws.Range("K136").Value = 68253468
$ws.Range("M136").Value = -68250918

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1866.7667
$ws.Range("I3").Value = 1536.3684
$ws.Range("J3").Value = 2437.4546
$ws.Range("K3").Value = 1536.3684
$ws.Range("L3").Value = 2437.4546
$ws.Range("M3").Value = -1422.3684
$ws.Range("N3").Value = -2665.4546
$ws.Range("H99").Value = 987.26666
$ws.Range("I99").Value = 848
$ws.Range("J99").Value = 1370.25
$ws.Range("K99").Value = 848
$ws.Range("L99").Value = 1370.25
$ws.Range("M99").Value = 650
$ws.Range("N99").Value = -4366.25
$ws.Range("H105").Value = 15153594
$ws.Range("I105").Value = 26317712
$ws.Range("K105").Value = 26317712
$ws.Range("M105").Value = -26315965
$ws.Range("H134").Value = 2842.75
$ws.Range("I134").Value = 1728.6428
$ws.Range("J134").Value = 4402.5
$ws.Range("K134").Value = 5185.928400000001
$ws.Range("L134").Value = 13207.5
$ws.Range("M134").Value = -2650.928400000001
$ws.Range("N134").Value = -18277.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4565.385
$ws.Range("I31").Value = 2591.6667
$ws.Range("J31").Value = 6257.143
$ws.Range("K31").Value = 2591.6667
$ws.Range("L31").Value = 6257.143
$ws.Range("M31").Value = -2296.6667
$ws.Range("N31").Value = -6847.143
$ws.Range("H34").Value = 4565.385
$ws.Range("I34").Value = 2591.6667
$ws.Range("J34").Value = 6257.143
$ws.Range("K34").Value = 2591.6667
$ws.Range("L34").Value = 6257.143
$ws.Range("M34").Value = -2389.6667
$ws.Range("N34").Value = -6661.143
$ws.Range("H132").Value = 36915.484
$ws.Range("I132").Value = 2425.0557
$ws.Range("J132").Value = 93354.37
$ws.Range("K132").Value = 7275.1671
$ws.Range("L132").Value = 280063.11
$ws.Range("M132").Value = -4745.1671
$ws.Range("N132").Value = -285123.11
$ws.Range("H134").Value = 29349.846
$ws.Range("I134").Value = 1717.1154
$ws.Range("J134").Value = 84615.30499999999
$ws.Range("K134").Value = 5151.3462
$ws.Range("L134").Value = 253845.915
$ws.Range("M134").Value = -2616.3462
$ws.Range("N134").Value = -258915.915

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 833.3333
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("H25").Value = 947.75
$ws.Range("J25").Value = 930.3333
$ws.Range("L25").Value = 2790.9999
$ws.Range("N25").Value = -3128.9999
$ws.Range("H30").Value = 947.75
$ws.Range("J30").Value = 930.3333
$ws.Range("L30").Value = 2790.9999
$ws.Range("N30").Value = -2994.9999
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -30588
$ws.Range("H80").Value = 1660
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1660
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4980
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -6852
$ws.Range("H83").Value = 1660
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1660
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 14940
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -24300
$ws.Range("H113").Value = 445.77274
$ws.Range("I113").Value = 329.91666
$ws.Range("K113").Value = 989.7499799999999
$ws.Range("M113").Value = 1180.25002
$ws.Range("H131").Value = 1204.3143
$ws.Range("I131").Value = 366.66666
$ws.Range("J131").Value = 1377.6207
$ws.Range("K131").Value = 1099.99998
$ws.Range("L131").Value = 4132.8621
$ws.Range("M131").Value = 3940.00002
$ws.Range("N131").Value = -14212.8621

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1020.61536
$ws.Range("I122").Value = 1000.6923
$ws.Range("J122").Value = 1040.5385
$ws.Range("K122").Value = 3002.0769
$ws.Range("L122").Value = 3121.6155
$ws.Range("M122").Value = -552.0769
$ws.Range("N122").Value = -8021.6155
$ws.Range("H132").Value = 51776.773
$ws.Range("I132").Value = 39973.348
$ws.Range("J132").Value = 73697.42999999999
$ws.Range("K132").Value = 119920.044
$ws.Range("L132").Value = 221092.29
$ws.Range("M132").Value = -117390.044
$ws.Range("N132").Value = -226152.29
$ws.Range("H134").Value = 34900
$ws.Range("J134").Value = 34900
$ws.Range("L134").Value = 104700
$ws.Range("N134").Value = -109770
$ws.Range("H135").Value = 34850.562
$ws.Range("J135").Value = 35126.668
$ws.Range("L135").Value = 35126.668
$ws.Range("N135").Value = -45266.668

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5428.2856
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 6166.3335
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 6166.3335
$ws.Range("M7").Value = -888
$ws.Range("N7").Value = -6390.3335
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null
$ws.Range("H122").Value = 2936.4849
$ws.Range("I122").Value = 2398
$ws.Range("J122").Value = 3582.6667
$ws.Range("K122").Value = 7194
$ws.Range("L122").Value = 10748.0001
$ws.Range("M122").Value = -4744
$ws.Range("N122").Value = -15648.0001
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H126").Value = 5428.2856
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 6166.3335
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 18499.0005
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -23439.0005
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1079.091
$ws.Range("I126").Value = 1087
$ws.Range("K126").Value = 3261
$ws.Range("M126").Value = -791
